$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for all data rows (2-115)
# from serial date 45175 to 45177.
$ws.Range("C2:C115").Value = 45177
